# Appends the Catholic calendar 'special days' entries for February
# through June (rows 12-39) to the Special_Days worksheet, mirroring
# the existing type/month/day/name_german layout used by rows 2-11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    ,@(2, 2, 'festivity', 'Darstellung des Herrn')
    ,@(2, 5, 'memory', 'Agatha')
    ,@(2, 6, 'memory', 'Paul Miki und Gefährten')
    ,@(2, 10, 'memory', 'Scholastika')
    ,@(2, 14, 'festivity', 'Cyrill und Methodius Slawenmission')
    ,@(2, 22, 'festivity', 'Kathedra Petri')
    ,@(2, 23, 'memory', 'Polykarp')
    ,@(2, 24, 'festivity', 'Matthias (Apstl.)')
    ,@(3, 7, 'memory', 'Perpetua und Felizitas')
    ,@(3, 19, 'high', 'Josef, Bräutigam der Gottesmutter Maria')
    ,@(3, 25, 'high', 'Verkündigung des Herrn')
    ,@(4, 7, 'memory', 'Johannes Baptist de la Salle')
    ,@(4, 11, 'memory', 'Stanislaus')
    ,@(4, 25, 'festivity', 'Markus (Evgl.)')
    ,@(4, 29, 'festivity', 'Katharina von Siena')
    ,@(5, 2, 'memory', 'Athanasius')
    ,@(5, 3, 'festivity', 'Philippus und Jakobus (Apstl.)')
    ,@(5, 21, 'memory', 'Christophorus Magallanes und Hermann Josef')
    ,@(5, 26, 'memory', 'Philipp Neri')
    ,@(6, 1, 'memory', 'Justin')
    ,@(6, 3, 'memory', 'Karl Lwanga und Gefährten')
    ,@(6, 5, 'festivity', 'Bonifatius')
    ,@(6, 11, 'memory', 'Barnabas (Apstl.)')
    ,@(6, 13, 'memory', 'Antonius von Padua')
    ,@(6, 21, 'memory', 'Aloisius Gonzaga')
    ,@(6, 24, 'high', 'Geburt Johannes des Täufers')
    ,@(6, 28, 'memory', 'Irenäus')
    ,@(6, 29, 'high', 'Petrus und Paulus (Apstl.)')
)

$startRow = 12
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $entry = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $entry[0]   # month
    $ws.Cells.Item($r, 2).Value = $entry[1]   # day
    $ws.Cells.Item($r, 3).Value = $entry[2]   # type
    $ws.Cells.Item($r, 4).Value = $entry[3]   # name_german
}
